$d = $word.ActiveDocument

# The primary header (header1.xml) holds most of the placeholder text that
# needs updating; the document body holds one more occurrence of "QWREW".
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)

# IMPORTANT ORDERING: "REW" is a substring of "QWREW", so the longer
# literal must be replaced first in every scope, otherwise the second
# pass would also clobber part of the text that used to read "QWREW".

# 1) Body: "QWREW" -> "QWR"  (case sensitive, exact match)
$d.Content.Find.Execute("QWREW", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "QWR", 2) | Out-Null

# 2) Header: "QWREW" -> "QWR"
$hdr.Range.Find.Execute("QWREW", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "QWR", 2) | Out-Null

# 3) Header: "REW" -> "QWER" (only the standalone occurrence remains now)
$hdr.Range.Find.Execute("REW", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "QWER", 2) | Out-Null

# 4) Header: "Rew" -> "Qwer" (5 occurrences)
$hdr.Range.Find.Execute("Rew", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Qwer", 2) | Out-Null

# 5) Header: "rew" -> "qwer" (3 occurrences)
$hdr.Range.Find.Execute("rew", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "qwer", 2) | Out-Null
